# Fix Trigonometry Lesson 15
# Adds a centered Law-of-Sines equation (sin a/a = sin b/b = sin g/c) as a
# new BodyText paragraph, right after the "All proportions will be equal."
# paragraph and before the "examples" bookmark / Examples heading.

$d = $word.ActiveDocument

# Locate the paragraph that ends the "All proportions will be equal." text.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*All proportions will be equal.*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'All proportions will be equal.' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# Insert a brand-new paragraph right after it (inherits the BodyText style
# from the source paragraph automatically).
[void]$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range

$mathXml = '<m:oMathPara xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><m:oMathParaPr><m:jc m:val="center"/></m:oMathParaPr><m:oMath><m:f><m:fPr><m:type m:val="bar"/></m:fPr><m:num><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>sin</m:t></m:r><m:r><m:t>α</m:t></m:r></m:num><m:den><m:r><m:t>a</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar"/></m:fPr><m:num><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>sin</m:t></m:r><m:r><m:t>β</m:t></m:r></m:num><m:den><m:r><m:t>b</m:t></m:r></m:den></m:f><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar"/></m:fPr><m:num><m:r><m:t>s</m:t></m:r><m:r><m:t>i</m:t></m:r><m:r><m:t>n</m:t></m:r><m:r><m:t>γ</m:t></m:r></m:num><m:den><m:r><m:t>c</m:t></m:r></m:den></m:f></m:oMath></m:oMathPara>'

# InsertXML replaces the new (empty) paragraph's contents with the equation
# and clears pPr in the process, so (re)apply the BodyText style afterwards.
[void]$newRange.InsertXML($mathXml)

$newPara2 = $d.Paragraphs.Item($targetIndex + 1)
$newPara2.Style = "BodyText"

Write-Host "Inserted Law of Sines equation paragraph after paragraph $targetIndex"
